$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")
$ws.Rows.Item(29).Insert()
$ws.Range("R29").Value = "lounge"
$ws.Range("S29").Value = "2024-09-03 13:08:08"
